# Add the "blackstone" (黑曜石) row to the MainIcon table.
#
# The table (表1) currently spans A1:N23 with header row 1 and data rows
# 2-23. We append one new data row (row 24) for the new "obsidian /
# blackstone" side-button entry, mirroring the previous row (SideButton2 /
# story) which shares the same D..M values (all zero except Flow=2) and
# the same ShowInScene/ShowInDungeon flags (false/true).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the Excel Table by one row - this grows ref/autoFilter (N23 -> N24)
# and the sheet dimension, and gives the new row the table's styling.
$tbl = $ws.ListObjects.Item(1)
$newListRow = $tbl.ListRows.Add()

# Seed row 24 from row 23: this clones number formats/styles (in particular
# the text-formatted TRUE/FALSE cells in columns G/H) so the new row matches
# the look of the existing rows exactly.
$ws.Range("A23:N23").Copy($ws.Range("A24:N24"))

# Now overwrite just the cells that differ for this new entry. Order matters
# for shared-string allocation: Icon (N) first, then Name (B), then
# Description (C), so new shared strings are appended as
# SideButton8, 黑曜石, 打开黑曜石面板 (matching authoring order).
$ws.Range("N24").Value = "SideButton8"
$ws.Range("B24").Value = "黑曜石"
$ws.Range("C24").Value = "打开黑曜石面板"
$ws.Range("A24").Value = 44

# Match the saved selection left by the author after adding the row.
$null = $ws.Range("C24").Select()
